$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$val) {
    # Force the value to be written as literal text, even if it looks numeric,
    # by building it via a formula and then converting the formula result to a
    # static value in place (keeps the original "General" cell style intact).
    $escaped = $val -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) | Out-Null
}
$excel.CutCopyMode = $false

# Row 2
Set-TextValue $ws.Range("D2") "51.041.29"

# Row 3
Set-TextValue $ws.Range("D3") "2.956.04"
$ws.Range("E3").Value = "  +0.43%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
Set-TextValue $ws.Range("D5") "380.65"
$ws.Range("E5").Value = "  +1.04%  "

# Row 6
Set-TextValue $ws.Range("D6") "102.03"
$ws.Range("E6").Value = "  -0.72%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.545"
$ws.Range("E7").Value = "  +1.70%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("E9").Value = "  +0.45%  "

# Row 10
Set-TextValue $ws.Range("D10") "36.53"
$ws.Range("E10").Value = "  -0.67%  "

# Row 11
$ws.Range("E11").Value = "  -0.91%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0852"
$ws.Range("E12").Value = "  +1.71%  "

# Row 13
Set-TextValue $ws.Range("D13") "3.421.69"
$ws.Range("E13").Value = "  +0.46%  "

# Row 14
$ws.Range("E14").Value = "  +2.20%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D15") "7.72"
$ws.Range("E15").Value = "  +5.04%  "

# Row 16
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D16") "12.08"
$ws.Range("E16").Value = "  +69.70%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.955.97"
$ws.Range("E17").Value = "  +0.42%  "

# Row 18
$ws.Range("E18").Value = "  +3.23%  "

# Row 19
Set-TextValue $ws.Range("D19") "51.103.50"
$ws.Range("E19").Value = "  -0.06%  "

# Row 20
$ws.Range("E20").Value = "  -2.45%  "

# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D21") "0.0₃0964"
$ws.Range("E21").Value = "  +0.85%  "

# Row 22
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D22") "12.36"
$ws.Range("E22").Value = "  -1.78%  "

# Row 23
$ws.Range("E23").Value = "  +16.05%  "

# Row 24
Set-TextValue $ws.Range("D24") "269.20"
$ws.Range("E24").Value = "  +2.24%  "

# Row 25
Set-TextValue $ws.Range("D25") "69.73"
$ws.Range("E25").Value = "  +2.18%  "

# Row 26
Set-TextValue $ws.Range("D26") "7.94"
$ws.Range("E26").Value = "  -2.58%  "

# Row 27
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("E28").Value = "  -1.06%  "

# Row 29
Set-TextValue $ws.Range("D29") "25.87"
$ws.Range("E29").Value = "  +0.73%  "

# Row 30
Set-TextValue $ws.Range("D30") "7.03"
$ws.Range("E30").Value = "  -10.68%  "

# Row 31
$ws.Range("E31").Value = "  -3.65%  "

# Row 32
Set-TextValue $ws.Range("D32") "10.42"
$ws.Range("E32").Value = "  +5.59%  "

# Row 33
$ws.Range("E33").Value = "  +5.73%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D34") "51.19"
$ws.Range("E34").Value = "  +0.65%  "

# Row 35
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D35") "34.31"
$ws.Range("E35").Value = "  +0.24%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.0435"
$ws.Range("E36").Value = "  -5.05%  "

# Row 37
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.27"
$ws.Range("E38").Value = "  +9.65%  "

# Row 39
$ws.Range("E39").Value = "  +1.91%  "

# Row 40
Set-TextValue $ws.Range("D40") "16.67"
$ws.Range("E40").Value = "  +1.17%  "

# Row 41
$ws.Range("E41").Value = "  +2.65%  "

# Row 42
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D42") "124.60"
$ws.Range("E42").Value = "  +2.25%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D43") "2.50"
$ws.Range("E43").Value = "  -3.34%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "21.69"
$ws.Range("E44").Value = "  +3.00%  "

# Row 45
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D45") "3.57"
$ws.Range("E45").Value = "  +10.56%  "

# Row 46
$ws.Range("E46").Value = "  +0.49%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.064.83"
$ws.Range("E47").Value = "  +3.26%  "

# Row 48
$ws.Range("E48").Value = "  -1.76%  "

# Row 49
$ws.Range("E49").Value = "  +1.93%  "

# Row 50
$ws.Range("E50").Value = "  -8.77%  "

# Row 51
Set-TextValue $ws.Range("D51") "5.37"
$ws.Range("E51").Value = "  +6.48%  "
